$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: date/time number format used by the "generate/handback" columns.
# ---------------------------------------------------------------------------
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ===========================================================================
# Sheet "Overview"  (new row 3: 9e8f5d6d-... report)
# ===========================================================================
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md"
$wsOv.Range("B3").Value = "e2e\9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-08-24 15:11:54"
$wsOv.Range("G3").NumberFormat = $dateFmt

# hyperlink for the new "Path And Name" cell, same repo/commit pattern as B2
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b06a634da378373e7f95371a64d877b4d013d53/e2e/9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md", `
    "", "", "e2e\9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md") | Out-Null

# ===========================================================================
# Sheet "zh-cn"  (row 2 refreshed with new handback run, new row 3 added)
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

# -- refresh existing row 2 with the new handoff/handback identifiers -------
$wsZh.Range("A2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.md"
$wsZh.Range("G2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.b90942244981dd365f980d378cfd591b010ac3f5.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-24 15:11:48"
$wsZh.Range("H2").NumberFormat = $dateFmt
$wsZh.Range("I2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.md"
$wsZh.Range("J2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.b90942244981dd365f980d378cfd591b010ac3f5.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-24 15:12:28"
$wsZh.Range("K2").NumberFormat = $dateFmt

$wsZh.Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b06a634da378373e7f95371a64d877b4d013d53/e2e/7a435067-3689-4895-975e-8d657d9fe8af.md", `
    "", "", "7a435067-3689-4895-975e-8d657d9fe8af.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e417ff402dd85df5f2ac160d676338b3878d796a/e2e/7a435067-3689-4895-975e-8d657d9fe8af.md", `
    "", "", "7a435067-3689-4895-975e-8d657d9fe8af.md") | Out-Null

# -- append the new row 3 for 9e8f5d6d... ------------------------------------
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.ca5c7e0a855523e7fe672fdee27784d59b2c309e.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-24 15:11:48"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("I3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md"
$wsZh.Range("J3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.ca5c7e0a855523e7fe672fdee27784d59b2c309e.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-24 15:12:28"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b06a634da378373e7f95371a64d877b4d013d53/e2e/9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md", `
    "", "", "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e417ff402dd85df5f2ac160d676338b3878d796a/e2e/9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md", `
    "", "", "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md") | Out-Null

# ===========================================================================
# Sheet "de-de"  (row 2 refreshed with new handback run, new row 3 added)
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

# -- refresh existing row 2 with the new handoff/handback identifiers -------
$wsDe.Range("A2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.md"
$wsDe.Range("G2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.b90942244981dd365f980d378cfd591b010ac3f5.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-24 15:11:54"
$wsDe.Range("H2").NumberFormat = $dateFmt
$wsDe.Range("I2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.md"
$wsDe.Range("J2").Value = "7a435067-3689-4895-975e-8d657d9fe8af.b90942244981dd365f980d378cfd591b010ac3f5.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-24 15:12:37"
$wsDe.Range("K2").NumberFormat = $dateFmt

$wsDe.Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b06a634da378373e7f95371a64d877b4d013d53/e2e/7a435067-3689-4895-975e-8d657d9fe8af.md", `
    "", "", "7a435067-3689-4895-975e-8d657d9fe8af.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1e88d5d1223df1dedeb4ca0b45bbbf2adfd23024/e2e/7a435067-3689-4895-975e-8d657d9fe8af.md", `
    "", "", "7a435067-3689-4895-975e-8d657d9fe8af.md") | Out-Null

# -- append the new row 3 for 9e8f5d6d... ------------------------------------
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.ca5c7e0a855523e7fe672fdee27784d59b2c309e.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-24 15:11:54"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("I3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md"
$wsDe.Range("J3").Value = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.ca5c7e0a855523e7fe672fdee27784d59b2c309e.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-24 15:12:37"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b06a634da378373e7f95371a64d877b4d013d53/e2e/9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md", `
    "", "", "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1e88d5d1223df1dedeb4ca0b45bbbf2adfd23024/e2e/9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md", `
    "", "", "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3.md") | Out-Null

Write-Output "Generate Report for Handback: done"
